# TestData.xlsx: trim the worksheet back down to an empty template.
#
# The "TestData" sheet used to ship with two fully-populated example rows
# (TC1/TC2) including mailto: hyperlinks in the username/password/email
# columns. Those sample values are removed so only the header row and the
# Test Case Name column remain populated; the now-unused strings fall out
# of the shared string table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Drop the mailto hyperlinks that lived on B2,B3,C2,C3,F2,F3 ...
[void]$ws.Range("B2:H3").Hyperlinks.Delete()
# ... and clear out the sample values in columns B:H for rows 2-3,
# keeping the Test Case Name values (column A) and cell styling intact.
[void]$ws.Range("B2:H3").ClearContents()

# Leave the sheet showing cell H14 selected, as it is in the edited file.
$ws.Activate()
[void]$ws.Range("H14").Select()
